$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain text in the source data (scraped),
# even when they look numeric (e.g. "6.990"). Force text format first so
# Excel does not reinterpret/renormalize them as numbers (dropping trailing
# zeros, stripping thousands separators, etc.).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.882.42"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.10"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.33"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3896"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.53"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.353"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08497"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.08"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.069"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.138"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.659.95"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.24"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07013"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.68"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.990"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.78"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.870.03"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.970"
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.12"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.11"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.444"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.05"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.892"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.496"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.834.44"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.018"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08214"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02914"
$ws.Range("E36").Value = "  -5.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.662"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.85"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2683"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09161"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7602"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.62"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6954"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.458"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08312"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.32"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.229"
$ws.Range("E51").Value = "  -2.98%  "
